$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.454.85'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.57%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.522.66'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.00%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.09'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.37%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('E8').Value = '  -2.53%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.523.60'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.94%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.96%  '

$ws.Range('E11').Value = '  -0.74%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.342'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.93%  '

$ws.Range('E13').Value = '  -3.03%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.984.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.02%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '70.229.70'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.84%  '

$ws.Range('E16').Value = '  -2.97%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.04'
$ws.Range('D17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.537.11'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.15%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.47'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.37%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.67'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.06%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '355.65'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.96%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.95'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.43%  '

$ws.Range('E23').Value = '  -3.28%  '

$ws.Range('E24').Value = '  +0.03%  '

$ws.Range('E25').Value = '  -4.02%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.08'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.11%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.55%  '

$ws.Range('E28').Value = '  -5.23%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.31%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0915'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.45%  '

$ws.Range('E31').Value = '  -2.86%  '

$ws.Range('E32').Value = '  +1.65%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '484.19'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.27%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.26%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.10%  '

$ws.Range('E36').Value = '  +5.23%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '156.92'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.68%  '

$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.65'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.34%  '

$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.90'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.16%  '

$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.65'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.43%  '

$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.321'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.25%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.74'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.92%  '

$ws.Range('E44').Value = '  -12.29%  '

$ws.Range('E45').Value = '  -7.10%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.33'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.87%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '143.23'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.84%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.54'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.36%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.528'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.11%  '

$ws.Range('E50').Value = '  -5.89%  '

$ws.Range('E51').Value = '  -0.76%  '
